$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meetups")

# Fix lab file names (renumbered to match the correct chapter/lab numbering)
$ws.Range("H5").Value = "/static/labs/05_foundations_for_inference.zip"
$ws.Range("H6").Value = "/static/labs/06_inference_for_categorial_data.zip"
$ws.Range("H7").Value = "/static/labs/07_inference_for_numerical_data"
$ws.Range("H8").Value = "/static/labs/08_simple_linear_regression.zip"
$ws.Range("H10").Value = "/static/labs/09_multiple_regression.zip"

# Add missing slides link for the Probability and Distributions meetup
$ws.Range("F4").Value = "/slides/03-Probability_and_Distributions.html"

# Update the active selection on the sheet
$ws.Range("F5").Select()
